$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Objetivos:" row (row 10) keeps its own cell references, but the text that
# used to live there ("Discutir as consequências...") was overwritten in the
# shared-strings table with the docente's name. Reproduce that by writing the
# new text straight into B10/C10.
$ws.Range("B10").Value = "230696 - Carlos José Todero Peixoto"
$ws.Range("C10").Value = "230696 - Carlos José Todero Peixoto"

# The old row 13 ("230696 - Carlos José Todero Peixoto", with no A label) is
# removed entirely; everything below shifts up by one row.
$ws.Rows("13").Delete()

# Former row 14 (now row 13, "Programa resumido:") gets new B/C content.
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Former row 16 (now row 15, "Programa:") gets new B/C content: the literal
# text "01/01/2018" (same string already used by B8/C8). Assigning it via
# .Value would be auto-coerced into a date serial number, so copy the
# existing text cell and paste values-only to keep it a plain text cell.
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("B8").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Former row 19 (now row 18, "Método:") gets new B/C content.
$ws.Range("B18").Value = "230696 - Carlos José Todero Peixoto"
$ws.Range("C18").Value = "230696 - Carlos José Todero Peixoto"

# Former row 20 (now row 19, "Critério:") gets new B/C content.
$ws.Range("B19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("C19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

# Former row 21 (now row 20, "Norma de recuperação:") gets new B/C content.
$ws.Range("B20").Value = "NF≥ 5,0."
$ws.Range("C20").Value = "NF≥ 5,0."

# Former row 22 (now row 21, "Bibliografia:") gets new B/C content.
$ws.Range("B21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("C21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
